$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5883377
$ws.Range("I106").Value = 7719871
$ws.Range("K106").Value = 7719871
$ws.Range("M106").Value = -7719240
$ws.Range("H113").Value = 16650.727
$ws.Range("I113").Value = 20154.8
$ws.Range("J113").Value = 13730.667
$ws.Range("K113").Value = 20154.8
$ws.Range("L113").Value = 13730.667
$ws.Range("M113").Value = -16900.8
$ws.Range("N113").Value = -20238.667
$ws.Range("H116").Value = 3672667.2
$ws.Range("I116").Value = 4405201
$ws.Range("K116").Value = 4405201
$ws.Range("M116").Value = -4401759
$ws.Range("H132").Value = 4083.9492
$ws.Range("I132").Value = 3595.25
$ws.Range("K132").Value = 10785.75
$ws.Range("M132").Value = -8255.75
$ws.Range("H133").Value = 143165.75
$ws.Range("J133").Value = 143165.75
$ws.Range("L133").Value = 143165.75
$ws.Range("N133").Value = -153285.75
$ws.Range("H137").Value = 8945.968999999999
$ws.Range("I137").Value = 10148.577
$ws.Range("K137").Value = 30445.731
$ws.Range("M137").Value = -27895.731
$ws.Range("H138").Value = 1771.5
$ws.Range("J138").Value = 3776.3333
$ws.Range("L138").Value = 11328.9999
$ws.Range("N138").Value = -21608.9999
$ws.Range("H141").Value = 3268.0952
$ws.Range("I141").Value = 3207.2354
$ws.Range("J141").Value = 3526.75
$ws.Range("K141").Value = 9621.706200000001
$ws.Range("L141").Value = 10580.25
$ws.Range("M141").Value = -4441.706200000001
$ws.Range("N141").Value = -20940.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 48800.59
$ws.Range("I2").Value = 3147.5386
$ws.Range("J2").Value = 114743.89
$ws.Range("K2").Value = 3147.5386
$ws.Range("L2").Value = 114743.89
$ws.Range("M2").Value = -3034.5386
$ws.Range("N2").Value = -114969.89
$ws.Range("H32").Value = 5382.8887
$ws.Range("I32").Value = 5486.4927
$ws.Range("K32").Value = 5486.4927
$ws.Range("M32").Value = -5199.4927
$ws.Range("H43").Value = 3500
$ws.Range("J43").Value = 3500
$ws.Range("L43").Value = 3500
$ws.Range("N43").Value = -4126
$ws.Range("H45").Value = 8181.2
$ws.Range("I45").Value = 7135.3335
$ws.Range("K45").Value = 7135.3335
$ws.Range("M45").Value = -6758.3335
$ws.Range("H61").Value = 4742.259
$ws.Range("I61").Value = 4751.9814
$ws.Range("K61").Value = 4751.9814
$ws.Range("M61").Value = -4539.9814
$ws.Range("H74").Value = 2835.2104
$ws.Range("I74").Value = 1101.8462
$ws.Range("K74").Value = 1101.8462
$ws.Range("M74").Value = -227.8462
$ws.Range("H77").Value = 2835.2104
$ws.Range("I77").Value = 1101.8462
$ws.Range("K77").Value = 5509.231
$ws.Range("M77").Value = -1141.231
$ws.Range("H110").Value = 2752.5417
$ws.Range("I110").Value = 1933.1428
$ws.Range("J110").Value = 3899.7
$ws.Range("K110").Value = 1933.1428
$ws.Range("L110").Value = 3899.7
$ws.Range("M110").Value = 111.8571999999999
$ws.Range("N110").Value = -7989.7
$ws.Range("H116").Value = 48800.59
$ws.Range("I116").Value = 3147.5386
$ws.Range("J116").Value = 114743.89
$ws.Range("K116").Value = 3147.5386
$ws.Range("L116").Value = 114743.89
$ws.Range("M116").Value = -853.5385999999999
$ws.Range("N116").Value = -119331.89
$ws.Range("H132").Value = 4081.9253
$ws.Range("I132").Value = 3696.0527
$ws.Range("K132").Value = 11088.1581
$ws.Range("M132").Value = -8558.158100000001
$ws.Range("H136").Value = 4742.259
$ws.Range("I136").Value = 4751.9814
$ws.Range("K136").Value = 14255.9442
$ws.Range("M136").Value = -11705.9442

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 48800.59
$ws.Range("I3").Value = 3147.5386
$ws.Range("J3").Value = 114743.89
$ws.Range("K3").Value = 3147.5386
$ws.Range("L3").Value = 114743.89
$ws.Range("M3").Value = -3033.5386
$ws.Range("N3").Value = -114971.89
$ws.Range("H95").Value = 45111.75
$ws.Range("J95").Value = 45111.75
$ws.Range("L95").Value = 45111.75
$ws.Range("N95").Value = -50603.75
$ws.Range("H107").Value = 3022.9167
$ws.Range("I107").Value = 2985.0588
$ws.Range("K107").Value = 2985.0588
$ws.Range("M107").Value = -1065.0588
$ws.Range("H134").Value = 6901.54
$ws.Range("I134").Value = 7023.933
$ws.Range("K134").Value = 21071.799
$ws.Range("M134").Value = -18536.799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4812203
$ws.Range("I99").Value = 9618210
$ws.Range("K99").Value = 9618210
$ws.Range("M99").Value = -9616712
$ws.Range("H126").Value = 4812203
$ws.Range("I126").Value = 9618210
$ws.Range("K126").Value = 28854630
$ws.Range("M126").Value = -28852160
$ws.Range("H132").Value = 18356
$ws.Range("I132").Value = 675.619
$ws.Range("K132").Value = 2026.857
$ws.Range("M132").Value = 503.143
$ws.Range("H134").Value = 1681.4
$ws.Range("I134").Value = 1611.25
$ws.Range("J134").Value = 1834.4546
$ws.Range("K134").Value = 4833.75
$ws.Range("L134").Value = 5503.3638
$ws.Range("M134").Value = -2298.75
$ws.Range("N134").Value = -10573.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 6990
$ws.Range("J106").Value = 6990
$ws.Range("L106").Value = 20970
$ws.Range("N106").Value = -22862
$ws.Range("H107").Value = 945.36
$ws.Range("J107").Value = 1146.5264
$ws.Range("L107").Value = 3439.5792
$ws.Range("N107").Value = -7279.5792
$ws.Range("H131").Value = 4753.5654
$ws.Range("J131").Value = 1977.7333
$ws.Range("L131").Value = 5933.199900000001
$ws.Range("N131").Value = -16013.1999
$ws.Range("H136").Value = 5171.8213
$ws.Range("I136").Value = 943.8095
$ws.Range("K136").Value = 2831.4285
$ws.Range("M136").Value = 2268.5715
$ws.Range("H138").Value = 805.8
$ws.Range("I138").Value = 805.8
$ws.Range("K138").Value = 2417.4
$ws.Range("M138").Value = 2722.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8249.65
$ws.Range("I70").Value = 8435.25
$ws.Range("K70").Value = 8435.25
$ws.Range("M70").Value = -8165.25
$ws.Range("H73").Value = 8249.65
$ws.Range("I73").Value = 8435.25
$ws.Range("K73").Value = 8435.25
$ws.Range("M73").Value = -7499.25
$ws.Range("H102").Value = 7607.6665
$ws.Range("I102").Value = 8129.2
$ws.Range("K102").Value = 8129.2
$ws.Range("M102").Value = -6507.2
$ws.Range("H123").Value = 43000
$ws.Range("J123").Value = 43000
$ws.Range("L123").Value = 43000
$ws.Range("N123").Value = -47900
$ws.Range("H132").Value = 1755.4333
$ws.Range("I132").Value = 1852.2858
$ws.Range("K132").Value = 5556.857400000001
$ws.Range("M132").Value = -3026.857400000001
$ws.Range("H139").Value = 42479
$ws.Range("J139").Value = 42479
$ws.Range("L139").Value = 42479
$ws.Range("N139").Value = -52759

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2848.3872
$ws.Range("I16").Value = 2936.96
$ws.Range("K16").Value = 2936.96
$ws.Range("M16").Value = -2766.96
$ws.Range("H61").Value = 4799.923
$ws.Range("I61").Value = 2957.0952
$ws.Range("J61").Value = 12539.8
$ws.Range("K61").Value = 2957.0952
$ws.Range("L61").Value = 12539.8
$ws.Range("M61").Value = -2755.0952
$ws.Range("N61").Value = -12943.8
$ws.Range("H82").Value = 2650.2
$ws.Range("I82").Value = 3595.25
$ws.Range("J82").Value = 1570.1428
$ws.Range("K82").Value = 3595.25
$ws.Range("L82").Value = 1570.1428
$ws.Range("M82").Value = -3234.25
$ws.Range("N82").Value = -2292.1428
$ws.Range("H85").Value = 2650.2
$ws.Range("I85").Value = 3595.25
$ws.Range("J85").Value = 1570.1428
$ws.Range("K85").Value = 3595.25
$ws.Range("L85").Value = 1570.1428
$ws.Range("M85").Value = -2347.25
$ws.Range("N85").Value = -4066.1428
$ws.Range("H113").Value = 4799.923
$ws.Range("I113").Value = 2957.0952
$ws.Range("J113").Value = 12539.8
$ws.Range("K113").Value = 2957.0952
$ws.Range("L113").Value = 12539.8
$ws.Range("M113").Value = -787.0952000000002
$ws.Range("N113").Value = -16879.8
$ws.Range("H136").Value = 3183.7778
$ws.Range("I136").Value = 1344.1936
$ws.Range("K136").Value = 4032.5808
$ws.Range("M136").Value = -1482.5808

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7590.4287
$ws.Range("J41").Value = 9650
$ws.Range("L41").Value = 9650
$ws.Range("N41").Value = -10430
$ws.Range("H96").Value = 12503419
$ws.Range("I96").Value = 16669934
$ws.Range("J96").Value = 3875
$ws.Range("K96").Value = 16669934
$ws.Range("L96").Value = 3875
$ws.Range("M96").Value = -16668561
$ws.Range("N96").Value = -6621
$ws.Range("H113").Value = 1871.6666
$ws.Range("I113").Value = 1063.2941
$ws.Range("J113").Value = 3834.8572
$ws.Range("K113").Value = 3189.8823
$ws.Range("L113").Value = 11504.5716
$ws.Range("M113").Value = -1019.8823
$ws.Range("N113").Value = -15844.5716
$ws.Range("H136").Value = 336048.97
$ws.Range("I136").Value = 453653.25
$ws.Range("K136").Value = 1360959.75
$ws.Range("M136").Value = -1358409.75
